$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Max/Min/Range helper block that lived in columns H:I (rows 2,3,5)
$ws.Range("H2:I2").ClearContents()
$ws.Range("H3:I3").ClearContents()
$ws.Range("H5:I5").ClearContents()

# H1 had the "Actual" label (shared string) with header styling; drop the text
# but keep the cell (and its style) in place.
$ws.Range("H1").ClearContents()

# Selection moves to H3 in the saved file
$ws.Range("H3").Select()
